$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Il7"
$ws.Cells.Item(2, 3).Value = "Il7r"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.8062
$ws.Cells.Item(2, 8).Value = 2.4186
$ws.Cells.Item(2, 9).Value = 0.581704519790725
$ws.Cells.Item(2, 10).Value = 0.5817045197907249
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.07349133333333334
$ws.Cells.Item(2, 14).Value = 0.220474
$ws.Cells.Item(2, 15).Value = 0.00285456188830886
$ws.Cells.Item(2, 16).Value = 0.00285456188830886
$ws.Cells.Item(2, 17).Value = 0.05924871293333334
$ws.Cells.Item(2, 18).Value = 0.5332384164
$ws.Cells.Item(2, 19).Value = 0.001660511552451611
$ws.Cells.Item(2, 20).Value = 0.00166051155245161

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Il7"
$ws.Cells.Item(3, 3).Value = "Il7r"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.8062
$ws.Cells.Item(3, 8).Value = 2.4186
$ws.Cells.Item(3, 9).Value = 0.581704519790725
$ws.Cells.Item(3, 10).Value = 0.5817045197907249
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.178101
$ws.Cells.Item(3, 14).Value = 0.534303
$ws.Cells.Item(3, 15).Value = 0.006917826957414881
$ws.Cells.Item(3, 16).Value = 0.006917826957414882
$ws.Cells.Item(3, 17).Value = 0.1435850262
$ws.Cells.Item(3, 18).Value = 1.2922652358
$ws.Cells.Item(3, 19).Value = 0.004024131208258356
$ws.Cells.Item(3, 20).Value = 0.004024131208258355

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Il7"
$ws.Cells.Item(4, 3).Value = "Il7r"
$ws.Cells.Item(4, 4).Value = "Resolving-Mac"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.8062
$ws.Cells.Item(4, 8).Value = 2.4186
$ws.Cells.Item(4, 9).Value = 0.581704519790725
$ws.Cells.Item(4, 10).Value = 0.5817045197907249
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 25.493631
$ws.Cells.Item(4, 14).Value = 76.48089300000001
$ws.Cells.Item(4, 15).Value = 0.9902276111542762
$ws.Cells.Item(4, 16).Value = 0.9902276111542762
$ws.Cells.Item(4, 17).Value = 20.5529653122
$ws.Cells.Item(4, 18).Value = 184.9766878098
$ws.Cells.Item(4, 19).Value = 0.576019877030015
$ws.Cells.Item(4, 20).Value = 0.5760198770300149

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Il7"
$ws.Cells.Item(5, 3).Value = "Il7r"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.4454536666666667
$ws.Cells.Item(5, 8).Value = 1.336361
$ws.Cells.Item(5, 9).Value = 0.3214120705251191
$ws.Cells.Item(5, 10).Value = 0.3214120705251191
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.07349133333333334
$ws.Cells.Item(5, 14).Value = 0.220474
$ws.Cells.Item(5, 15).Value = 0.00285456188830886
$ws.Cells.Item(5, 16).Value = 0.00285456188830886
$ws.Cells.Item(5, 17).Value = 0.03273698390155556
$ws.Cells.Item(5, 18).Value = 0.2946328551140001
$ws.Cells.Item(5, 19).Value = 0.0009174906469634445
$ws.Cells.Item(5, 20).Value = 0.0009174906469634445

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Il7"
$ws.Cells.Item(6, 3).Value = "Il7r"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.4454536666666667
$ws.Cells.Item(6, 8).Value = 1.336361
$ws.Cells.Item(6, 9).Value = 0.3214120705251191
$ws.Cells.Item(6, 10).Value = 0.3214120705251191
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.178101
$ws.Cells.Item(6, 14).Value = 0.534303
$ws.Cells.Item(6, 15).Value = 0.006917826957414881
$ws.Cells.Item(6, 16).Value = 0.006917826957414882
$ws.Cells.Item(6, 17).Value = 0.079335743487
$ws.Cells.Item(6, 18).Value = 0.714021691383
$ws.Cells.Item(6, 19).Value = 0.002223473085917202
$ws.Cells.Item(6, 20).Value = 0.002223473085917202

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Il7"
$ws.Cells.Item(7, 3).Value = "Il7r"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.4454536666666667
$ws.Cells.Item(7, 8).Value = 1.336361
$ws.Cells.Item(7, 9).Value = 0.3214120705251191
$ws.Cells.Item(7, 10).Value = 0.3214120705251191
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 25.493631
$ws.Cells.Item(7, 14).Value = 76.48089300000001
$ws.Cells.Item(7, 15).Value = 0.9902276111542762
$ws.Cells.Item(7, 16).Value = 0.9902276111542762
$ws.Cells.Item(7, 17).Value = 11.356231405597
$ws.Cells.Item(7, 18).Value = 102.206082650373
$ws.Cells.Item(7, 19).Value = 0.3182711067922384
$ws.Cells.Item(7, 20).Value = 0.3182711067922384

# Row 8
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Il7"
$ws.Cells.Item(8, 3).Value = "Il7r"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.1342733333333333
$ws.Cells.Item(8, 8).Value = 0.40282
$ws.Cells.Item(8, 9).Value = 0.09688340968415604
$ws.Cells.Item(8, 10).Value = 0.09688340968415604
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.07349133333333334
$ws.Cells.Item(8, 14).Value = 0.220474
$ws.Cells.Item(8, 15).Value = 0.00285456188830886
$ws.Cells.Item(8, 16).Value = 0.00285456188830886
$ws.Cells.Item(8, 17).Value = 0.009867926297777779
$ws.Cells.Item(8, 18).Value = 0.08881133667999999
$ws.Cells.Item(8, 19).Value = 0.0002765596888938053
$ws.Cells.Item(8, 20).Value = 0.0002765596888938053

# Row 9
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Il7"
$ws.Cells.Item(9, 3).Value = "Il7r"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.1342733333333333
$ws.Cells.Item(9, 8).Value = 0.40282
$ws.Cells.Item(9, 9).Value = 0.09688340968415604
$ws.Cells.Item(9, 10).Value = 0.09688340968415604
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.178101
$ws.Cells.Item(9, 14).Value = 0.534303
$ws.Cells.Item(9, 15).Value = 0.006917826957414881
$ws.Cells.Item(9, 16).Value = 0.006917826957414882
$ws.Cells.Item(9, 17).Value = 0.02391421494
$ws.Cells.Item(9, 18).Value = 0.21522793446
$ws.Cells.Item(9, 19).Value = 0.0006702226632393246
$ws.Cells.Item(9, 20).Value = 0.0006702226632393247

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Il7"
$ws.Cells.Item(10, 3).Value = "Il7r"
$ws.Cells.Item(10, 4).Value = "Resolving-Mac"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.1342733333333333
$ws.Cells.Item(10, 8).Value = 0.40282
$ws.Cells.Item(10, 9).Value = 0.09688340968415604
$ws.Cells.Item(10, 10).Value = 0.09688340968415604
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 25.493631
$ws.Cells.Item(10, 14).Value = 76.48089300000001
$ws.Cells.Item(10, 15).Value = 0.9902276111542762
$ws.Cells.Item(10, 16).Value = 0.9902276111542762
$ws.Cells.Item(10, 17).Value = 3.42311481314
$ws.Cells.Item(10, 18).Value = 30.80803331826
$ws.Cells.Item(10, 19).Value = 0.09593662733202291
$ws.Cells.Item(10, 20).Value = 0.09593662733202291
